$d = $word.ActiveDocument

# Locate the "Sim-real validated constants..." paragraph -- the two
# new paragraphs belong directly after it (and before "Usage rule...").
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "Sim-real validated constants in this run: WHEEL_RADIUS=0.049m, BASE_RADIUS=0.1085m, LIN_SCALE=1.0166, ANG_SCALE=1.2360, WZ_SIGN=-1.0.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph not found"
}

# Collapse to the end of the matched text (= end of that paragraph's
# run, right before its paragraph mark) and insert a fresh paragraph
# mark after it. InsertParagraphAfter does not advance the range's own
# Start/End, so remember the insertion point explicitly and step past
# the newly minted paragraph mark by one character to land inside it.
$anchor.Collapse(0)
$insertPos1 = $anchor.End
$anchor.InsertParagraphAfter()

$p1 = $d.Range($insertPos1 + 1, $insertPos1 + 1)
$p1.Text = "Supplementary check: sim_real_calibration_test.py validates pose-level distance/rotation only; deployment Go/NoGo is finalized by tune_sim_dynamics.py + replay_in_sim.py + check_calibration_gate.py."

# Repeat the same dance right after the paragraph we just filled in to
# add the second new paragraph.
$insertPos2 = $p1.End
$tail = $d.Range($insertPos2, $insertPos2)
$tail.InsertParagraphAfter()

$p2 = $d.Range($insertPos2 + 1, $insertPos2 + 1)
$p2.Text = "Synchronization: with AUTO_LOAD_COMP_FROM_DYNAMICS=True, sim_real_calibration_test.py first loads command_transform from calibration/tuned_dynamics.json and falls back to hardcoded constants if missing."

Write-Output "Inserted supplementary + synchronization paragraphs."
